$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the shared-string text from A2 up into A1, then drop the bold/border/
# center-aligned formatting that used to distinguish A1 (the old "0" header
# cell) from the plain data cell.
$ws.Range("A1").Value = 'questions = [
    {
        "title": "Jason enters a clothing store with the intention to steal a dress for his girlfriend\u2019s birthday. Jason grabs a dress and runs away without paying.What is the main legal issue in this legal scenario?",
        "ques_type": 2,
        "options": [
            "Whether Jason intended to steal the dress when he grabbed the dress.",
            "Whether Jason committed theft when he ran away with the dress without payment.",
            "Whether Jason entered the store with the intention to commit theft.",
            "Whether Jason committed burglary."
        ],
        "score": "Whether Jason committed theft when he ran away with the dress without payment."
    },
    {
        "title": "You are drafting an agreement for your client who acts as the authorized sales agent for Volt, a Memphis-based electric car company. Your client does not want Volt to appoint any other agent to advertise and sell Volt cars in Tennessee during the period of agreement. Volt agrees to this provision on the condition that Volt may market but not sell its products in Tennessee during the period of agreement.Which of the following clauses should you include in the agreement to reflect this understanding?",
        "ques_type": 2,
        "options": [
            "The Agent shall have sole advertising and selling rights subject to the company being permitted to market its products throughout Tennessee.",
            "The company shall not market and sell its products or appoint any other Agent to market and sell its products in Tennessee.",
            "The Agent shall have sole advertising and selling rights subject to the company being permitted to promote and sell its products throughout Tennessee.",
            "The company shall market its products but shall not appoint any other Agent to market and sell its products in Tennessee. "
        ],
        "score": "The Agent shall have sole advertising and selling rights subject to the company being permitted to market its products throughout Tennessee."
    },
    {
        "title": "You\u2019ve been given a list of the following sources as essential reading for research on WTO\u2019s Multi-Party Interim Appeal Arbitration Arrangement.Colombia \u2013 Anti-Dumping Duties on Frozen Fries from Belgium, Germany and the Netherlands WT/DS591/ARB25WTO Agreement: Marrakesh Agreement Establishing the World Trade OrganizationSungjoon Cho, \u201dA Global Constitutional Crisis\u201d  Florida State University Law Review, Vol. 49 issue 3 (Spring 2022)Lionel Bently, Brad Sherman, Dev Gangjee, and Phillip Johnson Intellectual Property Law (5th edition) Which of the above sources are secondary sources?",
        "ques_type": 2,
        "options": [
            "Only 1 and 2",
            "Only 2 and 4",
            "Only 2 and 3",
            "Only 3 and 4 "
        ],
        "score": "Only 3 and 4"
    },
    {
        "title": "You are a legal researcher working for a law firm specialized in criminal law. Your supervisor hands you a copy of a recent court case and asks you to read the judgment and write a basic summary of the court''s findings. Which of the following sentences uses the proper punctuation to indicate a direct quote from a judge in the decision?",
        "ques_type": 2,
        "options": [
            "The decision of the lower court is vacated-- stated Justice Warren --\"and the case is remanded for further proceedings.\"",
            "\"The decision of the lower court is vacated,\" stated Justice Warren, \"and the case is remanded for further proceedings.\"",
            "The decision of the lower court is vacated\" stated Justice Warren \"and the case is remanded for further proceedings.\"",
            "\"The decision of the lower court is vacated, stated Justice Warren and the case is remanded for further proceedings.\""
        ],
        "score": "\"The decision of the lower court is vacated,\" stated Justice Warren, \"and the case is remanded for further proceedings.\""
    }
]'
$ws.Range("A1").Style = "Normal"
$ws.Rows(1).AutoFit()

# Remove the now-redundant row that used to hold the shared string.
$ws.Rows(2).Delete()
